$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $val = $cell.Value2
    if ($null -eq $val -or $val -eq "") { continue }

    $parts = $val -split ", "
    $count = $parts.Count

    if ($count -lt 2) { continue }
    if ($parts[$count - 1] -ne "System") { continue }
    if ($parts -contains "admin@admin.com") { continue }

    # Remove the trailing "System" entry, then re-insert "System" right
    # after a leading lowercase "system" token if present, else at front.
    $rest = $parts[0..($count - 2)]

    if ($rest[0] -eq "system") {
        $newParts = @($rest[0], "System") + $rest[1..($rest.Count - 1)]
    } else {
        $newParts = @("System") + $rest
    }

    $cell.Value = ($newParts -join ", ")
}
